$d = $word.ActiveDocument

# Select "I have" (the part of the sentence being replaced) and type the
# replacement over it. Temporarily bolding the new text, then clearing the
# bold afterwards, prevents the engine from silently re-merging the newly
# typed run back into the untouched remainder of the sentence, matching the
# two separate <w:r> runs seen in the target diff (same rPr, split text).
$r = $d.Content
$r.Find.Execute("I have")
$r.Select()
$word.Selection.Font.Bold = 1
$word.Selection.TypeText("The cybersecurity analyst has")

$r2 = $d.Content
$r2.Find.Execute("The cybersecurity analyst has")
$r2.Font.Bold = 0
